$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 26406.678
$ws.Range("I98").Value = 1077.3334
$ws.Range("J98").Value = 79598.3
$ws.Range("K98").Value = 1077.3334
$ws.Range("L98").Value = 79598.3
$ws.Range("M98").Value = 420.6666
$ws.Range("N98").Value = -82594.3

$ws.Range("H113").Value = 2999.6667
$ws.Range("I113").Value = 2999
$ws.Range("K113").Value = 2999
$ws.Range("M113").Value = 255

$ws.Range("H122").Value = 26406.678
$ws.Range("I122").Value = 1077.3334
$ws.Range("J122").Value = 79598.3
$ws.Range("K122").Value = 3232.0002
$ws.Range("L122").Value = 238794.9
$ws.Range("M122").Value = -782.0001999999999
$ws.Range("N122").Value = -243694.9

$ws.Range("H129").Value = 402834.34
$ws.Range("J129").Value = 3856.0588
$ws.Range("L129").Value = 11568.1764
$ws.Range("N129").Value = -21568.1764

$ws.Range("H135").Value = 17242336
$ws.Range("I135").Value = 990.1786
$ws.Range("J135").Value = 500000000
$ws.Range("K135").Value = 8911.607399999999
$ws.Range("L135").Value = 4500000000
$ws.Range("M135").Value = -6376.607399999999
$ws.Range("N135").Value = -4500005070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 915
$ws.Range("I25").Value = 915
$ws.Range("K25").Value = 915
$ws.Range("M25").Value = -513

$ws.Range("H32").Value = 25587.684
$ws.Range("I32").Value = 24427.893
$ws.Range("K32").Value = 24427.893
$ws.Range("M32").Value = -24140.893

$ws.Range("H61").Value = 2426.6924
$ws.Range("I61").Value = 2261.2
$ws.Range("K61").Value = 2261.2
$ws.Range("M61").Value = -2049.2

$ws.Range("H80").Value = 58233
$ws.Range("J80").Value = 58233
$ws.Range("L80").Value = 58233
$ws.Range("N80").Value = -60229

$ws.Range("H83").Value = 58233
$ws.Range("J83").Value = 58233
$ws.Range("L83").Value = 174699
$ws.Range("N83").Value = -184683

$ws.Range("H122").Value = 2791.818
$ws.Range("I122").Value = 3487.2
$ws.Range("K122").Value = 10461.6
$ws.Range("M122").Value = -8011.599999999999

$ws.Range("H132").Value = 12197593
$ws.Range("I132").Value = 26317556
$ws.Range("J132").Value = 3079.9092
$ws.Range("K132").Value = 78952668
$ws.Range("L132").Value = 9239.7276
$ws.Range("M132").Value = -78950138
$ws.Range("N132").Value = -14299.7276

$ws.Range("H136").Value = 2426.6924
$ws.Range("I136").Value = 2261.2
$ws.Range("K136").Value = 6783.599999999999
$ws.Range("M136").Value = -4233.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 528.3333
$ws.Range("I64").Value = 342.5
$ws.Range("J64").Value = 900
$ws.Range("K64").Value = 342.5
$ws.Range("L64").Value = 900
$ws.Range("M64").Value = -117.5
$ws.Range("N64").Value = -1350

$ws.Range("H67").Value = 528.3333
$ws.Range("I67").Value = 342.5
$ws.Range("J67").Value = 900
$ws.Range("K67").Value = 342.5
$ws.Range("L67").Value = 900
$ws.Range("M67").Value = 437.5
$ws.Range("N67").Value = -2460

$ws.Range("H94").Value = 929.5
$ws.Range("I94").Value = 909
$ws.Range("J94").Value = 950
$ws.Range("K94").Value = 909
$ws.Range("L94").Value = 950
$ws.Range("M94").Value = -458
$ws.Range("N94").Value = -1852

$ws.Range("H130").Value = 48584.5
$ws.Range("J130").Value = 48584.5
$ws.Range("L130").Value = 48584.5
$ws.Range("N130").Value = -58624.5

$ws.Range("H132").Value = 50595.332
$ws.Range("J132").Value = 50595.332
$ws.Range("L132").Value = 50595.332
$ws.Range("N132").Value = -60715.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 46997.332
$ws.Range("J111").Value = 46997.332
$ws.Range("L111").Value = 46997.332
$ws.Range("N111").Value = -55177.332

$ws.Range("H122").Value = 55624036
$ws.Range("I122").Value = 71514984
$ws.Range("J122").Value = 5727.75
$ws.Range("K122").Value = 214544952
$ws.Range("L122").Value = 17183.25
$ws.Range("M122").Value = -214542502
$ws.Range("N122").Value = -22083.25

$ws.Range("H141").Value = 11348.857
$ws.Range("I141").Value = 13000
$ws.Range("J141").Value = 10898.546
$ws.Range("K141").Value = 13000
$ws.Range("L141").Value = 10898.546
$ws.Range("M141").Value = -7820
$ws.Range("N141").Value = -21258.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 465
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 763.63635
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 2290.90905
$ws.Range("M6").Value = -187
$ws.Range("N6").Value = -2516.90905

$ws.Range("H113").Value = 5289.1904
$ws.Range("I113").Value = 11620.333
$ws.Range("J113").Value = 540.8333
$ws.Range("K113").Value = 34860.999
$ws.Range("L113").Value = 1622.4999
$ws.Range("M113").Value = -32690.999
$ws.Range("N113").Value = -5962.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 38127090
$ws.Range("I20").Value = 152500000
$ws.Range("J20").Value = 2783.3333
$ws.Range("K20").Value = 152500000
$ws.Range("L20").Value = 2783.3333
$ws.Range("M20").Value = -152499755
$ws.Range("N20").Value = -3273.3333

$ws.Range("H24").Value = 34444850
$ws.Range("I24").Value = 103333336
$ws.Range("J24").Value = 602.8333
$ws.Range("K24").Value = 103333336
$ws.Range("L24").Value = 602.8333
$ws.Range("M24").Value = -103333163
$ws.Range("N24").Value = -948.8333

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H97").Value = 8101.35
$ws.Range("J97").Value = 13240.091
$ws.Range("L97").Value = 13240.091
$ws.Range("N97").Value = -14232.091

$ws.Range("H104").Value = 42165.75
$ws.Range("J104").Value = 42165.75
$ws.Range("L104").Value = 42165.75
$ws.Range("N104").Value = -49153.75

$ws.Range("H130").Value = 44183.5
$ws.Range("J130").Value = 44183.5
$ws.Range("L130").Value = 44183.5
$ws.Range("N130").Value = -54223.5

$ws.Range("H132").Value = 2263.7083
$ws.Range("I132").Value = 1858.3784
$ws.Range("J132").Value = 3627.0908
$ws.Range("K132").Value = 5575.135200000001
$ws.Range("L132").Value = 10881.2724
$ws.Range("M132").Value = -3045.135200000001
$ws.Range("N132").Value = -15941.2724

$ws.Range("H135").Value = 36043.945
$ws.Range("J135").Value = 36043.945
$ws.Range("L135").Value = 36043.945
$ws.Range("N135").Value = -46183.945

$ws.Range("H141").Value = 56139.2
$ws.Range("J141").Value = 56139.2
$ws.Range("L141").Value = 56139.2
$ws.Range("N141").Value = -66499.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 42420
$ws.Range("J121").Value = 42420
$ws.Range("L121").Value = 42420
$ws.Range("N121").Value = -45914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H132").Value = 1649.375
$ws.Range("I132").Value = 1208.5883
$ws.Range("J132").Value = 2719.8572
$ws.Range("K132").Value = 3625.7649
$ws.Range("L132").Value = 8159.571599999999
$ws.Range("M132").Value = -1095.7649
$ws.Range("N132").Value = -13219.5716

$ws.Range("H135").Value = 55466.445
$ws.Range("J135").Value = 55466.445
$ws.Range("L135").Value = 55466.445
$ws.Range("N135").Value = -65606.44500000001

$ws.Range("H137").Value = 63518
$ws.Range("J137").Value = 63518
$ws.Range("L137").Value = 63518
$ws.Range("N137").Value = -73718
